# Generate Report for Handoff
# The b.md file is now ready for handoff (a new handoff .xlf was produced),
# and its handback file turned out to be stale relative to the latest
# source, so the report rows for b.md move from
# "Handed back: in sync with en-US" to "Ready for handoff" and pick up the
# new handoff artifact name / timestamp plus an explanatory error message.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a56868816572aca1236a891e58f434d430ed0aa5/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fcad568bd17a7adad1079f009f5298bac53cbdc5/e2e/b.md."

# ---- Overview sheet : row 3 is the b.md file ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-18 08:39:01"

# ---- zh-cn sheet : row 3 is the b.md file ----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-18 08:38:55"
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Columns.Item(16).ColumnWidth = 39.17

# ---- de-de sheet : row 3 is the b.md file ----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-18 08:39:01"
$dede.Range("P3").Value = $errorDetail
$dede.Columns.Item(16).ColumnWidth = 39.17
